$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the Actual values for Sprint 2 (D3) and Sprint 3 (E3)
$ws.Range("D3").Value = 9
$ws.Range("E3").Value = 8

# Update the selection to match the recorded state
$ws.Range("C15").Select()
